$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C13").Value = 5
$ws.Range("C27").Value = 5
$ws.Range("E28").Value = 5

$ws.Application.ActiveWindow.ScrollRow = 8
$ws.Range("E28").Select()
